$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" column: header cell H1 (styled like the other header cells,
# e.g. G1 "sum") and data cell H2 with the save flag value.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
